# Applies the commit's change:
#  - add a new "Player Info" sheet as the first sheet, containing the
#    player's ID/NAME/BATTING_HAND/BOWL_STYLE
#  - keep the existing "ODI Batting" sheet as the second sheet, but
#    rewrite its MATCH_CARD_LINK column into a MATCH_CODE column that
#    only stores the numeric match code instead of the full URL

$wb = $excel.ActiveWorkbook

# The new sheet is inserted before the (current) active sheet, which
# puts it in position 1 and pushes "ODI Batting" to position 2.
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item(2)

# ---- Player Info sheet ----------------------------------------------
$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $playerInfo.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Keep "ID" as text ("4245") rather than a number, then drop back to the
# default "Normal" style so no explicit formatting lingers on the cell.
$idCell = $playerInfo.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "4245"
$idCell.Style = "Normal"

$playerInfo.Range("B2").Value = "Afsar Zazai"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

# ---- ODI Batting sheet: MATCH_CARD_LINK -> MATCH_CODE ---------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$usedRows = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Text
    if ($link -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.Style = "Normal"
    }
}
